# Update "paises" workbook:
#  - refresh the "last updated" timestamp
#  - update Noruega and Brasil stats
#  - Chile's numbers were updated (casos totales rose from 1306 to 1610),
#    which moves it up the (descending, sorted by column B) ranking from
#    row 32 to row 29; Luxemburgo, Ecuador and Japon each shift down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 15:12"

# --- Noruega (row 19): casos totales, nuevos casos, recuperados ---
$ws.Cells.Item(19, 2).Value = 3687
$ws.Cells.Item(19, 3).Value = 315
$ws.Cells.Item(19, 5).Value = 3665

# --- Brasil (row 24): casos totales, nuevos casos, recuperados ---
$ws.Cells.Item(24, 2).Value = 2991
$ws.Cells.Item(24, 3).Value = 6
$ws.Cells.Item(24, 5).Value = 2908

# --- Chile moves up to row 29 with its newly updated statistics, while
#     Luxemburgo, Ecuador and Japon shift down one row each (rows 30-32) ---

# Row 29: Chile (new figures)
$ws.Cells.Item(29, 1).Value = "Chile"
$ws.Cells.Item(29, 2).Value = 1610
$ws.Cells.Item(29, 3).Value = 304
$ws.Cells.Item(29, 4).Value = 22
$ws.Cells.Item(29, 5).Value = 1583
$ws.Cells.Item(29, 6).Value = 7
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 5

# Row 30: Luxemburgo (unchanged figures, shifted down from row 29)
$ws.Cells.Item(30, 1).Value = "Luxemburgo"
$ws.Cells.Item(30, 2).Value = 1453
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 6
$ws.Cells.Item(30, 5).Value = 1438
$ws.Cells.Item(30, 6).Value = 3
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 9

# Row 31: Ecuador (unchanged figures, shifted down from row 30)
$ws.Cells.Item(31, 1).Value = "Ecuador"
$ws.Cells.Item(31, 2).Value = 1403
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 3
$ws.Cells.Item(31, 5).Value = 1366
$ws.Cells.Item(31, 6).Value = 58
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 34

# Row 32: Japon (unchanged figures, shifted down from row 31)
$ws.Cells.Item(32, 1).Value = "Japon"
$ws.Cells.Item(32, 2).Value = 1387
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 359
$ws.Cells.Item(32, 5).Value = 981
$ws.Cells.Item(32, 6).Value = 57
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 47
